# ReceptiveField_Calculation.xlsx - "Add files via upload" edit
#
# 1. Fix the "rout" formula in the CNN table (Sheet2, column M, rows 16-24)
#    from (rin + kernel - 1) * jin  ->  rin + (kernel - 1) * jin
#    Each cell is edited individually so Excel keeps a standalone <f> per
#    cell (not a shared formula block), matching how the workbook looked
#    before the edit.
# 2. Move the active tab / selection from Sheet5!B3 to Sheet2!R8.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet2")

for ($r = 16; $r -le 24; $r++) {
    $ws1.Range("M$r").Formula = "=[@rin]+([@kernel]-1)*[@jin]"
}

$ws1.Activate()
$ws1.Range("R8").Select()
